$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit duplicates row 411 (Acelga / Vega Modelo de Temuco entry) by copying
# it and inserting the copy right below itself. Excel's "copy row, select the
# row below, Insert" performs an "Insert Copied Cells" - it pushes rows
# 412:465 down to 413:466 and fills the newly inserted row 412 with an exact
# copy of row 411 (values + formatting, including the date-formatted style on
# column D). This matches the diff: dimension grows from A1:R465 to A1:R466,
# row 411 is unchanged, and every row from 412 onward now holds what used to
# be one row above it, with old row 465's data ending up in new row 466.
$ws.Range("A411:R411").Copy()
$ws.Range("A412").EntireRow.Insert()
